# Horarios Línea 141 - refresh scrape (run at 02:42:39)
$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912": new arrival scraped, inserted as row 11 ---------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Update header/meta rows
$ws1.Range("A2").Value = "Última actualización: 02:42:39"
$ws1.Range("A3").Value = "Total filas: 8"

# Insert a new data row above the current row 11 (shifts old rows 11-12 down
# to 12-13) and fill it with the freshly scraped arrival.
$ws1.Rows.Item(11).Insert()
$ws1.Range("A11").Value = "02:42:39"
$ws1.Range("B11").Value = "03:54"
$ws1.Range("C11").Value = "14_ABASTO"
$ws1.Range("D11").Value = 72
$ws1.Range("E11").Value = "LP1912"

# --- Sheet "LP1912-215": only the "last updated" timestamp changes ---------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 02:42:39"

# --- Sheet "6203-6173": only the "last updated" timestamp changes ----------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 02:42:39"
